$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "307.27") as well as
# multi-dot "thousands" strings (e.g. "27.210.35") that are not valid numbers.
# Excel auto-converts plain numeric-looking text typed into a General cell into
# a real number (losing the original text + introducing float noise), so mark
# the whole price column as Text first, write the literal strings, then drop
# the temporary number format back to Normal so no stray style lingers.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = '27.210.35'
$ws.Range("E2").Value = '  +1.84%  '

$ws.Range("D3").Value = '1.906.84'
$ws.Range("E3").Value = '  +2.46%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.38%  '

$ws.Range("D5").Value = '307.27'
$ws.Range("E5").Value = '  +0.95%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").Value = '0.5225'
$ws.Range("E7").Value = '  +2.33%  '

$ws.Range("D8").Value = '0.3773'
$ws.Range("E8").Value = '  +3.39%  '

$ws.Range("D9").Value = '0.07252'
$ws.Range("E9").Value = '  +1.50%  '

$ws.Range("D10").Value = '21.25'
$ws.Range("E10").Value = '  +3.65%  '

$ws.Range("D11").Value = '0.8978'
$ws.Range("E11").Value = '  +1.02%  '

$ws.Range("D12").Value = '0.08351'
$ws.Range("E12").Value = '  +11.49%  '

$ws.Range("D13").Value = '1.915.24'
$ws.Range("E13").Value = '  +2.80%  '

$ws.Range("D14").Value = '94.86'
$ws.Range("E14").Value = '  +0.60%  '

$ws.Range("D15").Value = '5.269'
$ws.Range("E15").Value = '  +1.12%  '

$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  +0.38%  '

$ws.Range("D17").Value = '0.000008589'
$ws.Range("E17").Value = '  +3.18%  '

$ws.Range("D18").Value = '14.49'
$ws.Range("E18").Value = '  +2.58%  '

$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.26%  '

$ws.Range("D20").Value = '27.256.12'
$ws.Range("E20").Value = '  +1.84%  '

$ws.Range("D21").Value = '5.067'
$ws.Range("E21").Value = '  +1.66%  '

$ws.Range("D22").Value = '2.157.53'

$ws.Range("D23").Value = '10.62'
$ws.Range("E23").Value = '  +2.77%  '

$ws.Range("D24").Value = '6.446'
$ws.Range("E24").Value = '  +1.46%  '

$ws.Range("D25").Value = '2.289'
$ws.Range("E25").Value = '  +10.15%  '

$ws.Range("D26").Value = '145.94'
$ws.Range("E26").Value = '  +0.36%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '18.21'
$ws.Range("E27").Value = '  +2.16%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '1.743'
$ws.Range("E28").Value = '  -1.11%  '

$ws.Range("D29").Value = '114.74'
$ws.Range("E29").Value = '  +1.21%  '

$ws.Range("D30").Value = '4.974'
$ws.Range("E30").Value = '  +4.64%  '

$ws.Range("D31").Value = '4.798'
$ws.Range("E31").Value = '  +2.43%  '

$ws.Range("D32").Value = '0.09181'
$ws.Range("E32").Value = '  +0.55%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '0.8124'
$ws.Range("E33").Value = '  +9.32%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.05052'
$ws.Range("E34").Value = '  +0.41%  '

$ws.Range("D35").Value = '1.240'
$ws.Range("E35").Value = '  +7.53%  '

$ws.Range("D36").Value = '2.972'
$ws.Range("E36").Value = '  -0.55%  '

$ws.Range("D37").Value = '3.364'
$ws.Range("E37").Value = '  +4.84%  '

$ws.Range("D38").Value = '2.562'
$ws.Range("E38").Value = '  +2.86%  '

$ws.Range("D39").Value = '0.5715'
$ws.Range("E39").Value = '  +3.95%  '

$ws.Range("D40").Value = '0.01974'
$ws.Range("E40").Value = '  +0.46%  '

$ws.Range("D41").Value = '1.073'
$ws.Range("E41").Value = '  +0.68%  '

$ws.Range("D42").Value = '9.023'

$ws.Range("D43").Value = '6.581'
$ws.Range("E43").Value = '  +0.96%  '

$ws.Range("D44").Value = '117.92'
$ws.Range("E44").Value = '  +1.82%  '

$ws.Range("D45").Value = '0.1510'
$ws.Range("E45").Value = '  +2.04%  '

$ws.Range("D46").Value = '0.4831'
$ws.Range("E46").Value = '  +2.76%  '

$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  +0.36%  '

$ws.Range("D48").Value = '10.13'
$ws.Range("E48").Value = '  +1.55%  '

$ws.Range("D49").Value = '1.607'
$ws.Range("E49").Value = '  +3.53%  '

$ws.Range("D50").Value = '37.46'
$ws.Range("E50").Value = '  +2.15%  '

$ws.Range("D51").Value = '63.62'
$ws.Range("E51").Value = '  +1.44%  '

# Restore default styling on the price column now that the text values are set.
$priceCol.Style = "Normal"
